# The commit "Reads multiple excel sheets" renamed the single worksheet
# from "Sheet1" to "Sheet2" and left the cursor/selection on cell D15
# instead of E6 (as would happen after a user clicked around the sheet
# while wiring up the multi-sheet import code).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet: Sheet1 -> Sheet2
$ws.Name = "Sheet2"

# Move the active selection from E6 to D15
$ws.Range("D15").Select()
